$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 3 (Verde/1920 data), shifting it down to row 4.
$ws.Rows.Item(3).Insert()

# New row 3: updated data (Sin especificar, 600/1400/1500/1450)
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value = "Bíobío"
$ws.Cells.Item(3, 4).Value = 44489
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 300000000
$ws.Cells.Item(3, 7).Value = "Espárragos"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 600
$ws.Cells.Item(3, 11).Value = 1400
$ws.Cells.Item(3, 12).Value = 1500
$ws.Cells.Item(3, 13).Value = 1450
$ws.Cells.Item(3, 14).Value = "`$/kilo"
$ws.Cells.Item(3, 15).Value = "Provincia de Linares"
$ws.Cells.Item(3, 16).Value = 1450
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"
